$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 28.4088206944444
$ws.Range("D4").Value = 16.3245994675926
$ws.Range("B5").Value = 26.6388554143519
$ws.Range("D5").Value = 22.1135166608796
